# Generate Report for Handoff
# -----------------------------------------------------------------------
# This script rewrites the three worksheets ("Overview", "zh-cn", "de-de")
# of the localization-status workbook so that:
#   - the previously-failed handoff for 351eeb8b-....md is replaced by a
#     successful handoff for a new source file a7e92d36-....md (status
#     "Ready for handoff" instead of "Handoff transform failed"),
#   - a brand new source file ffff2779862e-....md shows up as a second
#     "Ready for handoff" row,
#   - the ".localization-config" bookkeeping row moves down to make room,
#   - the zh-cn / de-de detail sheets gain the corresponding "Latest
#     Handoff File" (.xlf) links + handoff datetimes and an "Include"
#     handoff reason for the two real files.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$mdCommit     = "e6d653be2f16ded94497fa795be84e444576a814"
$configCommit = "0eb6b76d2e933da6ea3fc8940a2b2ecaefe78969"

$urlNewMd      = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/a7e92d36-29de-4fd2-b571-1cfa984205c8.md"
$urlFfffMd     = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/ffff2779862e-1c90-401b-9786-9e4d2cf9062d.md"
$urlConfig     = "https://github.com/OpenLocalizationTest/oltest/blob/$configCommit/.localization-config"
$urlZhXlf      = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/a7e92d36-29de-4fd2-b571-1cfa984205c8.aa6171980a82ab00837cc07c12fda249f142ae56.zh-cn.xlf"
$urlDeXlf      = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/a7e92d36-29de-4fd2-b571-1cfa984205c8.aa6171980a82ab00837cc07c12fda249f142ae56.de-de.xlf"

$nameNewMd  = "a7e92d36-29de-4fd2-b571-1cfa984205c8.md"
$nameFfffMd = "ffff2779862e-1c90-401b-9786-9e4d2cf9062d.md"
$nameConfig = ".localization-config"
$nameZhXlf  = "a7e92d36-29de-4fd2-b571-1cfa984205c8.aa6171980a82ab00837cc07c12fda249f142ae56.zh-cn.xlf"
$nameDeXlf  = "a7e92d36-29de-4fd2-b571-1cfa984205c8.aa6171980a82ab00837cc07c12fda249f142ae56.de-de.xlf"

$readyForHandoff   = "Ready for handoff"
$notToBeLocalized  = "Not to be localized"
$zeroDate          = "0001-01-01 00:00:00"
$zhHandoffDate     = "2016-01-28 04:28:03"
$deHandoffDate     = "2016-01-28 04:28:13"

# =========================================================================
# Sheet "Overview" (File Name / zh-cn / de-de summary table)
# =========================================================================
$ovw = $wb.Worksheets.Item("Overview")

# Row 2: source file renamed + now ready for handoff (was failed)
$ovw.Range("A2").Value = $nameNewMd
$ovw.Range("B2").Value = $readyForHandoff
$ovw.Range("C2").Value = $readyForHandoff
$ovw.Hyperlinks.Add($ovw.Range("A2"), $urlNewMd, "", "", $nameNewMd) | Out-Null

# Row 3 (new): second source file, also ready for handoff
$ovw.Range("A3").Value = $nameFfffMd
$ovw.Range("B3").Value = $readyForHandoff
$ovw.Range("C3").Value = $readyForHandoff
$ovw.Hyperlinks.Add($ovw.Range("A3"), $urlFfffMd, "", "", $nameFfffMd) | Out-Null

# Row 4 (was row 3): .localization-config bookkeeping row, shifted down
$ovw.Range("A4").Value = $nameConfig
$ovw.Range("B4").Value = $notToBeLocalized
$ovw.Range("C4").Value = $notToBeLocalized
$ovw.Hyperlinks.Add($ovw.Range("A4"), $urlConfig, "", "", $nameConfig) | Out-Null

# =========================================================================
# Helper: populate one of the per-language detail sheets (zh-cn / de-de)
# =========================================================================
function Set-DetailSheet($ws, $xlfName, $xlfUrl, $handoffDate) {
    # Row 2: existing source file renamed, now included in a handoff
    $ws.Range("A2").Value = $nameNewMd
    $ws.Range("B2").Value = $readyForHandoff
    $ws.Range("C2").Value = $xlfName
    $ws.Range("D2").Value = $handoffDate
    $ws.Range("G2").Value = $zeroDate
    $ws.Range("H2").Value = "Include"
    $ws.Hyperlinks.Add($ws.Range("A2"), $urlNewMd, "", "", $nameNewMd) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $xlfUrl, "", "", $xlfName) | Out-Null

    # Row 3 (new): second source file, also included in a handoff
    $ws.Range("A3").Value = $nameFfffMd
    $ws.Range("B3").Value = $readyForHandoff
    $ws.Range("C3").Value = $xlfName
    $ws.Range("D3").Value = $handoffDate
    $ws.Range("G3").Value = $zeroDate
    $ws.Range("H3").Value = "Include"
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlFfffMd, "", "", $nameFfffMd) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $xlfUrl, "", "", $xlfName) | Out-Null

    # Row 4 (was row 3): .localization-config bookkeeping row, shifted down
    $ws.Range("A4").Value = $nameConfig
    $ws.Range("B4").Value = $notToBeLocalized
    $ws.Range("D4").Value = $zeroDate
    $ws.Range("G4").Value = $zeroDate
    $ws.Range("H4").Value = "Ignored"
    $ws.Hyperlinks.Add($ws.Range("A4"), $urlConfig, "", "", $nameConfig) | Out-Null
}

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$zh = $wb.Worksheets.Item("zh-cn")
Set-DetailSheet $zh $nameZhXlf $urlZhXlf $zhHandoffDate

# =========================================================================
# Sheet "de-de"
# =========================================================================
$de = $wb.Worksheets.Item("de-de")
Set-DetailSheet $de $nameDeXlf $urlDeXlf $deHandoffDate
